# Apply the diff: fix row 10 objective text, insert a "Docentes responsaveis" content row,
# shift rows 13-22 content down to 14-23 with corrected pairings, add bibliography text,
# add Requisitos row 23 and move the prerequisite text to row 24. Also narrow column A
# width declaration to col 1 only (col B keeps its own 60.71 width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (col A: min=1 max=1 now, instead of min=1 max=2) ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375
$ws.Columns.Item(2).ColumnWidth = 60.7109375
$ws.Columns.Item(3).ColumnWidth = 60.7109375

# --- Cell values (final, corrected layout) ---
$ws.Range('B1').Value = 'Ementa atual:'
$ws.Range('C1').Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range('B2').Value = 'LOQ4261'
$ws.Range('C2').Value = 'LOQ4261'
$ws.Range('A3').Value = 'Nome:'
$ws.Range('B3').Value = ' Planejamento, Programação e Controle da Produção'
$ws.Range('C3').Value = ' Planejamento, Programação e Controle da Produção'
$ws.Range('A4').Value = 'Name:'
$ws.Range('B4').Value = 'Production Planning, Scheduling and Control'
$ws.Range('C4').Value = 'Production Planning, Scheduling and Control'
$ws.Range('A5').Value = 'Créditos-aula:'
$ws.Range('B5').Value = '4'
$ws.Range('C5').Value = '4'
$ws.Range('A6').Value = 'Créditos-trabalho'
$ws.Range('B6').Value = '0'
$ws.Range('C6').Value = '0'
$ws.Range('A7').Value = 'Carga horária:'
$ws.Range('B7').Value = '60 h'
$ws.Range('C7').Value = '60 h'
$ws.Range('A8').Value = 'Ativação:'
$ws.Range('B8').Value = '01/01/2021'
$ws.Range('C8').Value = '01/01/2021'
$ws.Range('A9').Value = 'Semestre ideal:'
$ws.Range('B9').Value = 'EP-8'
$ws.Range('C9').Value = 'EP-8'
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = 'Apresentar um quadro conceitual de análise para auxiliar na formulação, avaliação e desenvolvimento de modelos para Planejamento, Programação e Controle da Produção nos diferentes ambientes de produção.'
$ws.Range('C10').Value = 'Apresentar um quadro conceitual de análise para auxiliar na formulação, avaliação e desenvolvimento de modelos para Planejamento, Programação e Controle da Produção nos diferentes ambientes de produção.'
$ws.Range('A11').Value = 'Objectives:'
$ws.Range('B11').Value = 'To present a conceptual framework of analysis to assist in the formulation, evaluation and development of models for Planning, Programming and Production Control in different production environments.'
$ws.Range('C11').Value = 'To present a conceptual framework of analysis to assist in the formulation, evaluation and development of models for Planning, Programming and Production Control in different production environments.'
$ws.Range('A12').Value = 'Docentes responsáveis:'
$ws.Range('B13').Value = '5701460 - Antonio Iacono'
$ws.Range('C13').Value = '5701460 - Antonio Iacono'
$ws.Range('A14').Value = 'Programa resumido:'
$ws.Range('B14').Value = '1. Caracterização do planejamento e controle da produção. 2. Gestão e previsão de demanda. 3. Planejamento agregado da produção. 4. Planejamento mestre da produção. 5. Planejamento e controle de estoques. 6. Planejamento de recursos de materiais (MRP). 7. Programação detalhada da produção. 8. Just In Time (JIT) e operações enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restrições (TOC). 11. Sistemas de controle da produção.'
$ws.Range('C14').Value = '1. Caracterização do planejamento e controle da produção. 2. Gestão e previsão de demanda. 3. Planejamento agregado da produção. 4. Planejamento mestre da produção. 5. Planejamento e controle de estoques. 6. Planejamento de recursos de materiais (MRP). 7. Programação detalhada da produção. 8. Just In Time (JIT) e operações enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restrições (TOC). 11. Sistemas de controle da produção.'
$ws.Range('A15').Value = 'Short syllabus:'
$ws.Range('B15').Value = '1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems.'
$ws.Range('C15').Value = '1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems.'
$ws.Range('A16').Value = 'Programa:'
$ws.Range('B16').Value = '1. Caracterização do planejamento e controle da produção. 2. Gestão e previsão de demanda. 3. Planejamento agregado da produção. 4. Planejamento mestre da produção. 5. Planejamento e controle de estoques. 6. Planejamento de recursos de materiais (MRP). 7. Programação detalhada da produção. 8. Just In Time (JIT) e operações enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restrições (TOC). 11. Sistemas de controle da produção.'
$ws.Range('C16').Value = '1. Caracterização do planejamento e controle da produção. 2. Gestão e previsão de demanda. 3. Planejamento agregado da produção. 4. Planejamento mestre da produção. 5. Planejamento e controle de estoques. 6. Planejamento de recursos de materiais (MRP). 7. Programação detalhada da produção. 8. Just In Time (JIT) e operações enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restrições (TOC). 11. Sistemas de controle da produção.'
$ws.Range('A17').Value = 'Syllabus:'
$ws.Range('B17').Value = '1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems.'
$ws.Range('C17').Value = '1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems.'
$ws.Range('A18').Value = 'Avaliação:'
$ws.Range('A19').Value = 'Método:'
$ws.Range('B19').Value = 'Aulas expositivas teóricas, aulas práticas, aulas de exercícios. MANTIDO'
$ws.Range('C19').Value = 'Aulas expositivas teóricas, aulas práticas, aulas de exercícios. MANTIDO'
$ws.Range('A20').Value = 'Critério:'
$ws.Range('B20').Value = 'M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.'
$ws.Range('C20').Value = 'M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.'
$ws.Range('A21').Value = 'Norma de recuperação:'
$ws.Range('B21').Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range('C21').Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range('A22').Value = 'Bibliografia:'
$ws.Range('B22').Value = '1.CORRÊA, H. L.; GIANESI, I. G. N.; CAON, M. Planejamento, programação e controle da produção: MRPII/ERP conceitos, uso e implantação. 5. ed. São Paulo: Atlas, 2007. 2.CORRÊA, H. L.; CORRÊA, C. A. Administração da Produção e Operações: manufatura e serviços: uma abordagem estratégica. 2.ed. São Paulo: Atlas, 2011. 3.DAVIS, M.M. et al. Fundamentos da administração da Produção. Porto Alegre: Bookman, 2018. 4.FERNANDES, F.C.F.; GODINHO FILHO. Planejamento e controle da produção: dos fundamentos ao essencial. São Paulo: Atlas, 2010. 5.GAITHER, N.; FRAZIER, G. Administração da Produção e Operações. 8. ed. São Paulo: Pioneira Thomson, 2005. 6.GONÇALVES, P.S. Administração de materiais. Rio de Janeiro: Elsevier, 2013. 7.HEIZER, J.; RENDER, B. Administração de Operações: bens e serviços. 5. ed. Rio de Janeiro: LTC, 2001.8.JACOBS, F.R.; CHASE, R. B. Administração da produção e de operações: o essencial. Porto Alegre: Bookman, 2009. 9.LUSTOSA, L. et. Al. Planejamento e controle da produção. Rio de Janeiro: Elsevier, 2008. 10.MOREIRA, D. A. Administração da Produção e Operações. 2. ed. São Paulo: Cengage Learning, 2008.11.REID, R.D.; SANDERS, N. R. Gestão de operações. Rio de Janeiro: LTC, 2005. 12.SLACK, N., BRANDON-JONES, A., JOHNSTON, R. Administração da produção. Henrique Luiz Corrêa (Trad.). 3. ed. São Paulo: Atlas, 2018. 13.TUBINO, D.V. Planejamento e controle da produção: teoria e prática. 2.ed. São Paulo: Atlas, 2009. 14.VOLLMANN, T.; BERRY, W.; WHYBARK, D.; JACOBS, F. Sistemas de planejamento e controle da produção: para o gerenciamento da cadeia de suprimentos. 5. ed. Porto Alegre: Bookman, 2006. 15.WANKE, P.F. Gerência de operações: uma abordagem logística. São Paulo: Atlas, 2010.'
$ws.Range('C22').Value = '1.CORRÊA, H. L.; GIANESI, I. G. N.; CAON, M. Planejamento, programação e controle da produção: MRPII/ERP conceitos, uso e implantação. 5. ed. São Paulo: Atlas, 2007. 2.CORRÊA, H. L.; CORRÊA, C. A. Administração da Produção e Operações: manufatura e serviços: uma abordagem estratégica. 2.ed. São Paulo: Atlas, 2011. 3.DAVIS, M.M. et al. Fundamentos da administração da Produção. Porto Alegre: Bookman, 2018. 4.FERNANDES, F.C.F.; GODINHO FILHO. Planejamento e controle da produção: dos fundamentos ao essencial. São Paulo: Atlas, 2010. 5.GAITHER, N.; FRAZIER, G. Administração da Produção e Operações. 8. ed. São Paulo: Pioneira Thomson, 2005. 6.GONÇALVES, P.S. Administração de materiais. Rio de Janeiro: Elsevier, 2013. 7.HEIZER, J.; RENDER, B. Administração de Operações: bens e serviços. 5. ed. Rio de Janeiro: LTC, 2001.8.JACOBS, F.R.; CHASE, R. B. Administração da produção e de operações: o essencial. Porto Alegre: Bookman, 2009. 9.LUSTOSA, L. et. Al. Planejamento e controle da produção. Rio de Janeiro: Elsevier, 2008. 10.MOREIRA, D. A. Administração da Produção e Operações. 2. ed. São Paulo: Cengage Learning, 2008.11.REID, R.D.; SANDERS, N. R. Gestão de operações. Rio de Janeiro: LTC, 2005. 12.SLACK, N., BRANDON-JONES, A., JOHNSTON, R. Administração da produção. Henrique Luiz Corrêa (Trad.). 3. ed. São Paulo: Atlas, 2018. 13.TUBINO, D.V. Planejamento e controle da produção: teoria e prática. 2.ed. São Paulo: Atlas, 2009. 14.VOLLMANN, T.; BERRY, W.; WHYBARK, D.; JACOBS, F. Sistemas de planejamento e controle da produção: para o gerenciamento da cadeia de suprimentos. 5. ed. Porto Alegre: Bookman, 2006. 15.WANKE, P.F. Gerência de operações: uma abordagem logística. São Paulo: Atlas, 2010.'
$ws.Range('A23').Value = 'Requisitos:'
$ws.Range('B24').Value = 'LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)
'
$ws.Range('C24').Value = 'LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)
'

# --- Clear cells that are no longer used at their old positions ---
$ws.Range('A13').Value = $null
$ws.Range('B18').Value = $null
$ws.Range('C18').Value = $null
$ws.Range('B23').Value = $null
$ws.Range('C23').Value = $null

# --- Reset row heights that must go back to default (no explicit height) ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(23).AutoFit()

# --- Apply explicit custom row heights for the final layout ---
$ws.Rows.Item(10).RowHeight = 60.0
$ws.Rows.Item(11).RowHeight = 60.0
$ws.Rows.Item(14).RowHeight = 60.0
$ws.Rows.Item(15).RowHeight = 60.0
$ws.Rows.Item(16).RowHeight = 120.0
$ws.Rows.Item(17).RowHeight = 120.0
$ws.Rows.Item(19).RowHeight = 60.0
$ws.Rows.Item(20).RowHeight = 60.0
$ws.Rows.Item(21).RowHeight = 60.0
$ws.Rows.Item(22).RowHeight = 120.0
$ws.Rows.Item(24).RowHeight = 30.0
